$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.521.28'
$ws.Range("E2").Value = '  -1.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.277.93'
$ws.Range("E3").Value = '  +0.88%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '94.60'
$ws.Range("E5").Value = '  -4.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '266.39'
$ws.Range("E6").Value = '  -2.48%  '
$ws.Range("E7").Value = '  -1.07%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -3.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '44.38'
$ws.Range("E10").Value = '  -7.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0932'
$ws.Range("E11").Value = '  -1.32%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.70'
$ws.Range("E12").Value = '  -6.31%  '
$ws.Range("E13").Value = '  +0.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.619.55'
$ws.Range("E14").Value = '  +1.01%  '
$ws.Range("E15").Value = '  -2.24%  '
$ws.Range("E16").Value = '  +1.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.281.60'
$ws.Range("E17").Value = '  +1.50%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.479.16'
$ws.Range("E18").Value = '  -1.45%  '
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.17'
$ws.Range("E20").Value = '  -0.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.96'
$ws.Range("E21").Value = '  +1.56%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.36'
$ws.Range("E22").Value = '  -0.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.99'
$ws.Range("E23").Value = '  -0.72%  '
$ws.Range("E24").Value = '  -8.86%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("E26").Value = '  -0.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.20'
$ws.Range("E27").Value = '  -2.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.47'
$ws.Range("E28").Value = '  -1.48%  '
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.95'
$ws.Range("E30").Value = '  -2.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.17'
$ws.Range("E31").Value = '  +0.85%  '
$ws.Range("E32").Value = '  +2.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0880'
$ws.Range("E33").Value = '  -3.71%  '
$ws.Range("E34").Value = '  -6.16%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("E37").Value = '  -5.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.40'
$ws.Range("E38").Value = '  +0.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.29'
$ws.Range("E39").Value = '  -7.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.32'
$ws.Range("E40").Value = '  +5.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.234'
$ws.Range("E41").Value = '  -6.31%  '
$ws.Range("E42").Value = '  +14.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.86'
$ws.Range("E43").Value = '  -4.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.76'
$ws.Range("E44").Value = '  +1.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.78'
$ws.Range("E45").Value = '  +3.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.20'
$ws.Range("E46").Value = '  -4.94%  '
$ws.Range("E47").Value = '  -1.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '97.34'
$ws.Range("E48").Value = '  -3.20%  '
$ws.Range("E49").Value = '  -1.00%  '
$ws.Range("E50").Value = '  +3.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.499.32'
$ws.Range("E51").Value = '  +1.11%  '
